# Update the weekly Fruta/Hortaliza prices ("Fruta / hortaliza, semanal").
# Each row's Fecha/Volumen/Precio values shift to the following week's
# figures (row 5 wraps around to the original row 2 figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44708
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("S2").Value = 1025

$ws.Range("D3").Value = 44357
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("S3").Value = 725

$ws.Range("D4").Value = 44533
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16500
$ws.Range("S4").Value = 825

$ws.Range("D5").Value = 44320
$ws.Range("M5").Value = 80
